$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ADDS Chart")
Write-Host $ws.Name
